# Append 4 new daily rows (2025-11-14 .. 2025-11-17) to the "Chart" sheet,
# and update the "Videos" rollup on the "Table" sheet to match.

$wb = $excel.ActiveWorkbook
$chart = $wb.Worksheets.Item("Chart")
$table = $wb.Worksheets.Item("Table")

# New data rows for the Chart sheet: Date, No video indexed, Video indexed, Impressions
$newRows = @(
    @{ Row = 42; Date = "2025-11-14"; NoVideoIndexed = 23; VideoIndexed = 1; Impressions = 0 },
    @{ Row = 43; Date = "2025-11-15"; NoVideoIndexed = 23; VideoIndexed = 1; Impressions = 0 },
    @{ Row = 44; Date = "2025-11-16"; NoVideoIndexed = 23; VideoIndexed = 1; Impressions = 0 },
    @{ Row = 45; Date = "2025-11-17"; NoVideoIndexed = 23; VideoIndexed = 1; Impressions = "" }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    # Column A holds a date-formatted string ("2025-11-14"). Writing a
    # date-looking string straight into a General cell gets auto-parsed
    # into a date serial by the COM layer, so briefly mark the cell as
    # Text first, then strip the explicit formatting back off afterwards
    # (ClearFormats keeps the already-committed string value) so the
    # stored cell ends up with the default style, same as its neighbors.
    $cellA = $chart.Cells.Item($rowNum, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $r.Date
    $cellA.ClearFormats()

    $chart.Cells.Item($rowNum, 2).Value = $r.NoVideoIndexed
    $chart.Cells.Item($rowNum, 3).Value = $r.VideoIndexed

    $cellD = $chart.Cells.Item($rowNum, 4)
    if ($r.Impressions -eq "") {
        # Last row's Impressions is blank text rather than 0.
        $cellD.NumberFormat = "@"
        $cellD.Value = " "
        $cellD.Value = ""
        $cellD.ClearFormats()
    } else {
        $cellD.Value = $r.Impressions
    }
}

# "Table" sheet: the "Videos" count for the existing validation row drops
# from 24 to 23 now that one more video has been indexed.
$table.Range("C2").Value = 23
